$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next 4 reporting dates to column A (rows 120-123)
$ws.Range("A120").Value = "29 05 2020"
$ws.Range("A121").Value = "30 05 2020"
$ws.Range("A122").Value = "31 05 2020"
$ws.Range("A123").Value = "01 06 2020"

# Per-state numeric data for the newly populated rows 117-120
$newData = @{
    "B117" = 34.35316206269
    "C117" = 31.461567974975
    "D117" = 9.682141444889799
    "F117" = 24.534816570541
    "G117" = 9.9042134163728
    "H117" = 18.009538566191
    "I117" = 17.744587076551
    "J117" = 56.261395292848
    "K117" = 17.819177815061
    "L117" = 14.581659545417
    "M117" = 15.154890628045
    "O117" = 9.5212119488586
    "P117" = 20.672724956881
    "Q117" = 13.315999747757
    "R117" = 18.783205340643
    "S117" = 18.238562052807
    "T117" = 16.067581625616
    "U117" = 24.824284513937
    "V117" = 13.462587754797
    "W117" = 24.292405097361
    "X117" = 21.162607165411
    "Y117" = 26.304868232186
    "Z117" = 19.126088876036
    "AA117" = 17.993302111515
    "AB117" = 20.104919832136
    "AD117" = 5.8353767234635
    "AE117" = 3.9099752861643
    "AF117" = 16.071638986805
    "AG117" = 48.992419650377
    "AH117" = 8.011057911075101
    "AI117" = 76.464547603555
    "AJ117" = 27.095108544429
    "AK117" = 22.992724522718
    "AL117" = 25.824251930266
    "AM117" = 17.258483083738
    "AN117" = 19.95497599782
    "AO117" = 16.141133450342
    "AP117" = 3.7387583781419
    "AQ117" = 14.313975784886
    "AS117" = 23.520662390165
    "AT117" = 5.2651456417518
    "AU117" = 0
    "AV117" = 12.629980107353
    "AW117" = 12.567674037473
    "AX117" = 18.204474174019
    "AY117" = 32.310924613671
    "BA117" = 3.6719786796652
    "BB117" = 14.950642786179
    "BC117" = 15.435777693133
    "BD117" = 10.481469747349
    "BE117" = 18.784880893357
    "B118" = 26.697730813912
    "C118" = 33.306399061421
    "D118" = 7.7716216819366
    "F118" = 22.295785257544
    "G118" = 9.8609753062695
    "H118" = 14.54587099109
    "I118" = 21.121348619232
    "J118" = 56.357269963873
    "K118" = 12.208265793873
    "L118" = 15.321996807082
    "M118" = 15.461881412475
    "O118" = 18.16039432267
    "P118" = 18.784440737015
    "Q118" = 27.305048192037
    "R118" = 19.125364529924
    "S118" = 20.071860176876
    "T118" = 16.65668158002
    "U118" = 25.012834387445
    "V118" = 13.281739661556
    "W118" = 20.681082230809
    "X118" = 21.425731666838
    "Y118" = 22.029343732872
    "Z118" = 17.732842381566
    "AA118" = 20.718213769303
    "AB118" = 16.107489291758
    "AD118" = 4.3882280913761
    "AE118" = 28.867418248908
    "AF118" = 17.616336043588
    "AG118" = 100.6034170318
    "AH118" = 10.570936229946
    "AI118" = 60.36136425668
    "AJ118" = 25.061183410066
    "AK118" = 22.259395334311
    "AL118" = 20.378164648417
    "AM118" = 17.632696474895
    "AN118" = 19.022125963684
    "AO118" = 11.575885287139
    "AP118" = 3.7373750871337
    "AQ118" = 16.588643730519
    "AS118" = 16.992048949138
    "AT118" = 5.065297295052
    "AU118" = 17.080443781854
    "AV118" = 11.194987896274
    "AW118" = 12.648418323732
    "AX118" = 25.757899524294
    "AY118" = 31.200521639934
    "BA118" = 0.92795837399062
    "BB118" = 10.819756225473
    "BC118" = 10.773371180161
    "BD118" = 16.681983930871
    "BE118" = 15.033540794491
    "B119" = 19.612749745187
    "C119" = 28.904345372195
    "D119" = 5.9367745895823
    "F119" = 20.212330491732
    "G119" = 10.453637582865
    "H119" = 11.616673445176
    "I119" = 26.593621105162
    "J119" = 43.363817463451
    "K119" = 7.4897432503356
    "L119" = 12.759614544252
    "M119" = 16.342859892792
    "O119" = 15.712730221456
    "P119" = 12.905006946854
    "Q119" = 31.334366505694
    "R119" = 20.225299074036
    "S119" = 15.262442717391
    "T119" = 22.290371342047
    "U119" = 21.879859510669
    "V119" = 9.874547665176999
    "W119" = 20.423037610349
    "X119" = 16.653120337102
    "Y119" = 17.743743018998
    "Z119" = 15.823293374388
    "AA119" = 19.184288230055
    "AB119" = 19.628299671174
    "AD119" = 3.0974662759549
    "AE119" = 22.957216037338
    "AF119" = 18.518185575456
    "AG119" = 84.672264231139
    "AH119" = 13.264068425655
    "AI119" = 45.349873687044
    "AJ119" = 22.429948709845
    "AK119" = 15.105904713753
    "AL119" = 15.285332987769
    "AM119" = 16.682381245011
    "AN119" = 17.116048460989
    "AO119" = 11.666229712074
    "AP119" = 4.156740225037
    "AQ119" = 14.379695826089
    "AS119" = 11.426645824794
    "AT119" = 10.450062819342
    "AU119" = 15.18748826046
    "AV119" = 21.369165750207
    "AW119" = 12.150974068978
    "AX119" = 24.538158898366
    "AY119" = 26.956024057795
    "BA119" = 19.404937588851
    "BB119" = 10.248978662861
    "BC119" = 15.008425431878
    "BD119" = 21.926636173541
    "BE119" = 11.414943856816
    "B120" = 13.32365269742
    "C120" = 24.820140947707
    "D120" = 4.2486678180018
    "F120" = 16.528969778114
    "G120" = 10.140283217175
    "H120" = 9.2742608546016
    "I120" = 20.318777969026
    "J120" = 31.533400730774
    "K120" = 3.6749308181075
    "L120" = 10.986699158903
    "M120" = 18.832241511016
    "O120" = 23.60499955144
    "P120" = 12.547994565143
    "Q120" = 26.020612961396
    "R120" = 17.469544848248
    "S120" = 13.087251045331
    "T120" = 22.425815553931
    "U120" = 19.007614023498
    "V120" = 6.8691130201028
    "W120" = 23.435733571108
    "X120" = 21.392083943849
    "Y120" = 13.619147874391
    "Z120" = 16.091685477995
    "AA120" = 13.73666222201
    "AB120" = 17.526552193139
    "AD120" = 8.295849785462799
    "AE120" = 17.572512308571
    "AF120" = 19.354560055537
    "AG120" = 68.438928797703
    "AH120" = 9.205855068886899
    "AI120" = 40.539021503028
    "AJ120" = 17.41684418515
    "AK120" = 9.05722870072
    "AL120" = 18.821304125523
    "AM120" = 16.387996652959
    "AN120" = 15.381877191516
    "AO120" = 7.931894438355
    "AP120" = 2.0672388973494
    "AQ120" = 18.43958511277
    "AS120" = 6.8099260846574
    "AT120" = 15.519300294733
    "AU120" = 33.354746776157
    "AV120" = 24.957865762314
    "AW120" = 12.415777399898
    "AX120" = 19.281041794621
    "AY120" = 21.982460696777
    "BA120" = 35.127483307948
    "BB120" = 13.108066451965
    "BC120" = 15.119539674056
    "BD120" = 17.028707977123
    "BE120" = 8.062473863743101
}

foreach ($addr in $newData.Keys) {
    $ws.Range($addr).Value = $newData[$addr]
}

Write-Output "Updated rows 117-123 with new data"
